$d = $word.ActiveDocument

$d.Content.Find.Execute("Start Fight", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start fight", 2)

$d.Content.Find.Execute("Pick up item.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Pick up item", 2)

$d.Content.Find.Execute("Start new single-player", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start new single-player game", 2)

$d.Content.Find.Execute("Start multiplayer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start multiplayer game", 2)
